$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New donation rows to append below the existing data (rows 2-3 already exist)
$data = @(
    @(45211, "Victory Drones",     350000),
    @(45211, "General donations",  3000),
    @(45194, "General donations",  15000),
    @(45192, "General donations",  3000),
    @(45191, "General donations",  3089.78),
    @(45191, "General donations",  12000),
    @(45189, "General donations",  3000),
    @(45163, "General donations",  10000),
    @(45161, "General donations",  3500),
    @(45160, "General donations",  10000),
    @(45152, "General donations",  250000),
    @(45148, "General donations",  840000),
    @(45098, "General donations",  3456),
    @(45068, "General donations",  3673)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $entry[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
